# "More graphs wrt CN"
# - Rename "Sheet3" to "errors" and populate it with a new error/404 reference table.
# - Move the active/selected tab from "example table" to the new "errors" sheet.
# - Update the remembered selection on "example table" from J39 to J30.

$wb = $excel.ActiveWorkbook

$wsSummary = $wb.Worksheets.Item("Sheet1")
$wsData    = $wb.Worksheets.Item("example table")
$wsNotes   = $wb.Worksheets.Item("Sheet3")

# --- Re-sequence the merged A-column header cells on Sheet1 ---------------------
# (cosmetic re-serialisation order only; ranges/formatting are unchanged)
$mergedOrder = @("A22:A25","A2:A5","A6:A9","A10:A13","A14:A17","A18:A21","A26:A29","A30:A33","A34:A37","A38:A41","A42:A45")
foreach ($ref in $mergedOrder) {
    $wsSummary.Range($ref).UnMerge() | Out-Null
}
foreach ($ref in $mergedOrder) {
    $wsSummary.Range($ref).Merge() | Out-Null
}

# Rename the (previously empty) third sheet to its new purpose.
$wsNotes.Name = "errors"

# --- Fill in the "errors" sheet -------------------------------------------------
# Written in the same order the original author entered it (this drives shared
# string allocation order, not just final cell content).
$wsNotes.Range("A1").Value = "Voor 404's"
$wsNotes.Range("A2").Value = "scriptname"
$wsNotes.Range("A3").Value = "date_cet"
$wsNotes.Range("A5").Value = "page_seq"
$wsNotes.Range("A6").Value = "topdomain"
$wsNotes.Range("A7").Value = "urlnoparams"
$wsNotes.Range("A8").Value = "status_code"
$wsNotes.Range("A9").Value = "error_code"
$wsNotes.Range("A4").Value = "task_succeed"
$wsNotes.Range("A10").Value = "number"
$wsNotes.Range("A11").Value = "total_time_sec"
$wsNotes.Range("A12").Value = "sec_per_page"
$wsNotes.Range("B12").Value = "Influence on page loading time"
$wsNotes.Range("B4").Value = "only task_succeed = 1 has influence on reported page loading times."

$wsNotes.Columns.Item(1).ColumnWidth = 14.8

# --- Update view/selection state -------------------------------------------------
# Move the "example table" selection before handing activation to "errors".
$wsData.Activate() | Out-Null
$wsData.Range("J30").Select() | Out-Null

$wsNotes.Activate() | Out-Null
$wsNotes.Range("B5").Select() | Out-Null
try { $excel.ActiveWindow.ScrollRow = 17 } catch {}

Write-Output "done"
